# Auto-generated Excel COM-interop script to apply the Gilgamesh_Profits.xlsx diff.
# Updates computed currentAveragePrice / LevePrice / LeveProfit figures across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets (columns H-N) per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80,8).Value = 504.7619  # H80
$ws.Cells.Item(80,10).Value = 425.57144  # J80
$ws.Cells.Item(80,12).Value = 1276.71432  # L80
$ws.Cells.Item(80,14).Value = -3272.71432  # N80
$ws.Cells.Item(83,8).Value = 504.7619  # H83
$ws.Cells.Item(83,10).Value = 425.57144  # J83
$ws.Cells.Item(83,12).Value = 3830.14296  # L83
$ws.Cells.Item(83,14).Value = -13814.14296  # N83
$ws.Cells.Item(86,8).Value = 62502388  # H86
$ws.Cells.Item(86,10).Value = 3100.5  # J86
$ws.Cells.Item(86,12).Value = 3100.5  # L86
$ws.Cells.Item(86,14).Value = -5346.5  # N86
$ws.Cells.Item(88,8).Value = 6734065  # H88
$ws.Cells.Item(88,9).Value = 25250350  # I88
$ws.Cells.Item(88,10).Value = 870.4545000000001  # J88
$ws.Cells.Item(88,11).Value = 25250350  # K88
$ws.Cells.Item(88,12).Value = 870.4545000000001  # L88
$ws.Cells.Item(88,13).Value = -25249944  # M88
$ws.Cells.Item(88,14).Value = -1682.4545  # N88
$ws.Cells.Item(89,8).Value = 62502388  # H89
$ws.Cells.Item(89,10).Value = 3100.5  # J89
$ws.Cells.Item(89,12).Value = 15502.5  # L89
$ws.Cells.Item(89,14).Value = -26734.5  # N89
$ws.Cells.Item(91,8).Value = 6734065  # H91
$ws.Cells.Item(91,9).Value = 25250350  # I91
$ws.Cells.Item(91,10).Value = 870.4545000000001  # J91
$ws.Cells.Item(91,11).Value = 25250350  # K91
$ws.Cells.Item(91,12).Value = 870.4545000000001  # L91
$ws.Cells.Item(91,13).Value = -25248946  # M91
$ws.Cells.Item(91,14).Value = -3678.4545  # N91
$ws.Cells.Item(112,8).Value = 2022.2667  # H112
$ws.Cells.Item(112,10).Value = 2113.1428  # J112
$ws.Cells.Item(112,12).Value = 6339.428400000001  # L112
$ws.Cells.Item(112,14).Value = -8555.428400000001  # N112
$ws.Cells.Item(132,8).Value = 6573.2085  # H132
$ws.Cells.Item(132,9).Value = 6785.087  # I132
$ws.Cells.Item(132,10).Value = 1700  # J132
$ws.Cells.Item(132,11).Value = 20355.261  # K132
$ws.Cells.Item(132,12).Value = 5100  # L132
$ws.Cells.Item(132,13).Value = -17825.261  # M132
$ws.Cells.Item(132,14).Value = -10160  # N132
$ws.Cells.Item(137,8).Value = 2554.2856  # H137
$ws.Cells.Item(137,9).Value = 1739.4242  # I137
$ws.Cells.Item(137,10).Value = 15999.5  # J137
$ws.Cells.Item(137,11).Value = 5218.2726  # K137
$ws.Cells.Item(137,12).Value = 47998.5  # L137
$ws.Cells.Item(137,13).Value = -2668.2726  # M137
$ws.Cells.Item(137,14).Value = -53098.5  # N137

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 5343.8984  # H32
$ws.Cells.Item(32,9).Value = 4457.4155  # I32
$ws.Cells.Item(32,11).Value = 4457.4155  # K32
$ws.Cells.Item(32,13).Value = -4170.4155  # M32
$ws.Cells.Item(97,8).Value = 1516.7222  # H97
$ws.Cells.Item(97,9).Value = 1619.0714  # I97
$ws.Cells.Item(97,10).Value = 1158.5  # J97
$ws.Cells.Item(97,11).Value = 1619.0714  # K97
$ws.Cells.Item(97,12).Value = 1158.5  # L97
$ws.Cells.Item(97,13).Value = -1123.0714  # M97
$ws.Cells.Item(97,14).Value = -2150.5  # N97
$ws.Cells.Item(102,8).Value = 3997.4412  # H102
$ws.Cells.Item(102,9).Value = 3921.606  # I102
$ws.Cells.Item(102,11).Value = 3921.606  # K102
$ws.Cells.Item(102,13).Value = -2299.606  # M102
$ws.Cells.Item(110,8).Value = 2013.2903  # H110
$ws.Cells.Item(110,9).Value = 929.0909  # I110
$ws.Cells.Item(110,11).Value = 929.0909  # K110
$ws.Cells.Item(110,13).Value = 1115.9091  # M110

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20,8).Value = 12502169  # H20
$ws.Cells.Item(20,10).Value = 1076.8182  # J20
$ws.Cells.Item(20,12).Value = 1076.8182  # L20
$ws.Cells.Item(20,14).Value = -1570.8182  # N20
$ws.Cells.Item(60,8).Value = 0  # H60
$ws.Cells.Item(60,10).Value = 0  # J60
$ws.Cells.Item(60,12).Value = 0  # L60
$ws.Cells.Item(60,14).ClearContents()  # N60
$ws.Cells.Item(94,8).Value = 38462204  # H94
$ws.Cells.Item(94,9).Value = 50000460  # I94
$ws.Cells.Item(94,11).Value = 50000460  # K94
$ws.Cells.Item(94,13).Value = -50000009  # M94
$ws.Cells.Item(99,8).Value = 129685.5  # H99
$ws.Cells.Item(99,9).Value = 252502.5  # I99
$ws.Cells.Item(99,11).Value = 252502.5  # K99
$ws.Cells.Item(99,13).Value = -251004.5  # M99
$ws.Cells.Item(105,8).Value = 7431591.5  # H105
$ws.Cells.Item(105,9).Value = 457485.88  # I105
$ws.Cells.Item(105,10).Value = 19233924  # J105
$ws.Cells.Item(105,11).Value = 457485.88  # K105
$ws.Cells.Item(105,12).Value = 19233924  # L105
$ws.Cells.Item(105,13).Value = -455738.88  # M105
$ws.Cells.Item(105,14).Value = -19237418  # N105
$ws.Cells.Item(134,8).Value = 2314.3333  # H134
$ws.Cells.Item(134,9).Value = 1736.0741  # I134
$ws.Cells.Item(134,11).Value = 5208.2223  # K134
$ws.Cells.Item(134,13).Value = -2673.2223  # M134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4,8).Value = 0  # H4
$ws.Cells.Item(4,9).Value = 0  # I4
$ws.Cells.Item(4,11).Value = 0  # K4
$ws.Cells.Item(4,13).ClearContents()  # M4
$ws.Cells.Item(31,8).Value = 3495.52  # H31
$ws.Cells.Item(31,9).Value = 2718.2593  # I31
$ws.Cells.Item(31,11).Value = 2718.2593  # K31
$ws.Cells.Item(31,13).Value = -2423.2593  # M31
$ws.Cells.Item(34,8).Value = 3495.52  # H34
$ws.Cells.Item(34,9).Value = 2718.2593  # I34
$ws.Cells.Item(34,11).Value = 2718.2593  # K34
$ws.Cells.Item(34,13).Value = -2516.2593  # M34
$ws.Cells.Item(105,8).Value = 1993.9412  # H105
$ws.Cells.Item(105,9).Value = 1095.2858  # I105
$ws.Cells.Item(105,11).Value = 1095.2858  # K105
$ws.Cells.Item(105,13).Value = 651.7141999999999  # M105
$ws.Cells.Item(107,8).Value = 769.8182  # H107
$ws.Cells.Item(107,9).Value = 790.8889  # I107
$ws.Cells.Item(107,10).Value = 675  # J107
$ws.Cells.Item(107,11).Value = 790.8889  # K107
$ws.Cells.Item(107,12).Value = 675  # L107
$ws.Cells.Item(107,13).Value = 1129.1111  # M107
$ws.Cells.Item(107,14).Value = -4515  # N107
$ws.Cells.Item(141,8).Value = 344583.25  # H141
$ws.Cells.Item(141,10).Value = 344583.25  # J141
$ws.Cells.Item(141,12).Value = 344583.25  # L141
$ws.Cells.Item(141,14).Value = -354943.25  # N141

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3,8).Value = 7872.5557  # H3
$ws.Cells.Item(3,9).Value = 7836.75  # I3
$ws.Cells.Item(3,11).Value = 23510.25  # K3
$ws.Cells.Item(3,13).Value = -23398.25  # M3
$ws.Cells.Item(69,8).Value = 1574.8  # H69
$ws.Cells.Item(69,10).Value = 2324.6667  # J69
$ws.Cells.Item(69,12).Value = 6974.000100000001  # L69
$ws.Cells.Item(69,14).Value = -8596.000100000001  # N69
$ws.Cells.Item(72,8).Value = 1574.8  # H72
$ws.Cells.Item(72,10).Value = 2324.6667  # J72
$ws.Cells.Item(72,12).Value = 20922.0003  # L72
$ws.Cells.Item(72,14).Value = -29034.0003  # N72
$ws.Cells.Item(122,8).Value = 2339.3076  # H122
$ws.Cells.Item(122,10).Value = 2401.0908  # J122
$ws.Cells.Item(122,12).Value = 21609.8172  # L122
$ws.Cells.Item(122,14).Value = -26509.8172  # N122

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5,8).Value = 10000  # H5
$ws.Cells.Item(5,10).Value = 0  # J5
$ws.Cells.Item(5,12).Value = 0  # L5
$ws.Cells.Item(5,14).ClearContents()  # N5
$ws.Cells.Item(24,8).Value = 15041.625  # H24
$ws.Cells.Item(24,9).Value = 10999  # I24
$ws.Cells.Item(24,10).Value = 15217.392  # J24
$ws.Cells.Item(24,11).Value = 10999  # K24
$ws.Cells.Item(24,12).Value = 15217.392  # L24
$ws.Cells.Item(24,13).Value = -10826  # M24
$ws.Cells.Item(24,14).Value = -15563.392  # N24
$ws.Cells.Item(80,8).Value = 250004240  # H80
$ws.Cells.Item(80,9).Value = 333336830  # I80
$ws.Cells.Item(80,10).Value = 6500  # J80
$ws.Cells.Item(80,11).Value = 333336830  # K80
$ws.Cells.Item(80,12).Value = 6500  # L80
$ws.Cells.Item(80,13).Value = -333335832  # M80
$ws.Cells.Item(80,14).Value = -8496  # N80
$ws.Cells.Item(83,8).Value = 250004240  # H83
$ws.Cells.Item(83,9).Value = 333336830  # I83
$ws.Cells.Item(83,10).Value = 6500  # J83
$ws.Cells.Item(83,11).Value = 1666684150  # K83
$ws.Cells.Item(83,12).Value = 32500  # L83
$ws.Cells.Item(83,13).Value = -1666679158  # M83
$ws.Cells.Item(83,14).Value = -42484  # N83
$ws.Cells.Item(97,8).Value = 506.9091  # H97
$ws.Cells.Item(97,9).Value = 673.4  # I97
$ws.Cells.Item(97,10).Value = 368.16666  # J97
$ws.Cells.Item(97,11).Value = 673.4  # K97
$ws.Cells.Item(97,12).Value = 368.16666  # L97
$ws.Cells.Item(97,13).Value = -177.4  # M97
$ws.Cells.Item(97,14).Value = -1360.16666  # N97
$ws.Cells.Item(106,8).Value = 30000  # H106
$ws.Cells.Item(106,10).Value = 30000  # J106
$ws.Cells.Item(106,12).Value = 30000  # L106
$ws.Cells.Item(106,14).Value = -32524  # N106
$ws.Cells.Item(122,8).Value = 1772.1052  # H122
$ws.Cells.Item(122,9).Value = 1748.4445  # I122
$ws.Cells.Item(122,11).Value = 5245.333500000001  # K122
$ws.Cells.Item(122,13).Value = -2795.333500000001  # M122
$ws.Cells.Item(132,8).Value = 3439.1025  # H132
$ws.Cells.Item(132,9).Value = 2681.6775  # I132
$ws.Cells.Item(132,10).Value = 6374.125  # J132
$ws.Cells.Item(132,11).Value = 8045.032499999999  # K132
$ws.Cells.Item(132,12).Value = 19122.375  # L132
$ws.Cells.Item(132,13).Value = -5515.032499999999  # M132
$ws.Cells.Item(132,14).Value = -24182.375  # N132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 6955.933  # H7
$ws.Cells.Item(7,9).Value = 4760.4287  # I7
$ws.Cells.Item(7,10).Value = 8877  # J7
$ws.Cells.Item(7,11).Value = 4760.4287  # K7
$ws.Cells.Item(7,12).Value = 8877  # L7
$ws.Cells.Item(7,13).Value = -4648.4287  # M7
$ws.Cells.Item(7,14).Value = -9101  # N7
$ws.Cells.Item(16,8).Value = 1063.2142  # H16
$ws.Cells.Item(16,9).Value = 1063.2142  # I16
$ws.Cells.Item(16,11).Value = 1063.2142  # K16
$ws.Cells.Item(16,13).Value = -893.2141999999999  # M16
$ws.Cells.Item(40,8).Value = 24404.715  # H40
$ws.Cells.Item(40,9).Value = 31868.945  # I40
$ws.Cells.Item(40,10).Value = 3734.5386  # J40
$ws.Cells.Item(40,11).Value = 31868.945  # K40
$ws.Cells.Item(40,12).Value = 3734.5386  # L40
$ws.Cells.Item(40,13).Value = -31732.945  # M40
$ws.Cells.Item(40,14).Value = -4006.5386  # N40
$ws.Cells.Item(100,8).Value = 3449.3333  # H100
$ws.Cells.Item(100,9).Value = 3422.15  # I100
$ws.Cells.Item(100,10).Value = 3585.25  # J100
$ws.Cells.Item(100,11).Value = 3422.15  # K100
$ws.Cells.Item(100,12).Value = 3585.25  # L100
$ws.Cells.Item(100,13).Value = -2881.15  # M100
$ws.Cells.Item(100,14).Value = -4667.25  # N100
$ws.Cells.Item(122,8).Value = 4613.5713  # H122
$ws.Cells.Item(122,9).Value = 4389.1816  # I122
$ws.Cells.Item(122,10).Value = 5436.3335  # J122
$ws.Cells.Item(122,11).Value = 13167.5448  # K122
$ws.Cells.Item(122,12).Value = 16309.0005  # L122
$ws.Cells.Item(122,13).Value = -10717.5448  # M122
$ws.Cells.Item(122,14).Value = -21209.0005  # N122
$ws.Cells.Item(126,8).Value = 6955.933  # H126
$ws.Cells.Item(126,9).Value = 4760.4287  # I126
$ws.Cells.Item(126,10).Value = 8877  # J126
$ws.Cells.Item(126,11).Value = 14281.2861  # K126
$ws.Cells.Item(126,12).Value = 26631  # L126
$ws.Cells.Item(126,13).Value = -11811.2861  # M126
$ws.Cells.Item(126,14).Value = -31571  # N126

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2,8).Value = 1667583.4  # H2
$ws.Cells.Item(2,9).Value = 1667583.4  # I2
$ws.Cells.Item(2,11).Value = 1667583.4  # K2
$ws.Cells.Item(2,13).Value = -1667471.4  # M2
$ws.Cells.Item(107,8).Value = 425.56522  # H107
$ws.Cells.Item(107,9).Value = 387.95  # I107
$ws.Cells.Item(107,11).Value = 1163.85  # K107
$ws.Cells.Item(107,13).Value = 756.1500000000001  # M107
$ws.Cells.Item(122,8).Value = 14709317  # H122
$ws.Cells.Item(122,9).Value = 3144.7693  # I122
$ws.Cells.Item(122,11).Value = 9434.3079  # K122
$ws.Cells.Item(122,13).Value = -6984.3079  # M122
$ws.Cells.Item(132,8).Value = 3559.5688  # H132
$ws.Cells.Item(132,9).Value = 3505.6956  # I132
$ws.Cells.Item(132,11).Value = 10517.0868  # K132
$ws.Cells.Item(132,13).Value = -7987.086800000001  # M132

